$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 (shifts existing rows 78..178 down to 79..179,
# duplicating the formatting of the row above as Excel normally does on insert).
$ws.Rows(78).Insert()

# Populate the newly inserted row 78 with the new data point.
$ws.Cells.Item(78, 1).Value = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44664
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 100112037
$ws.Cells.Item(78, 7).Value = "Cebollín"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 1600
$ws.Cells.Item(78, 11).Value = 1100
$ws.Cells.Item(78, 12).Value = 1200
$ws.Cells.Item(78, 13).Value = 1150
$ws.Cells.Item(78, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(78, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(78, 16).Value = 192
$ws.Cells.Item(78, 17).Value = 6
$ws.Cells.Item(78, 18).Value = "Hortaliza"
